$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 295
$ws1.Range("F4").Value = 1145
$ws1.Range("F5").Value = 590

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 8

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 295
$ws4.Range("F4").Value = 1145
$ws4.Range("F5").Value = 8
$ws4.Range("F6").Value = 590
